$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9794444441795349
$ws.Range("B1").Value = 1.411014437675476
$ws.Range("C1").Value = 2.55925989151001
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.591899991035461
